$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'328.96"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.22%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'41.20"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'4.16%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.633"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.11%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.08201"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'2.48%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'8.755"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'1.88%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'2.002"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.87%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'4.486"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.17%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-0.95%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9200"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-0.61%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1277"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'3.02%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1952"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-1.22%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.09280"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'1.24%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03894"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'6.97%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.1059"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.91%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001310"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'1.27%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.006150"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'0.23%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E19").Value = "'2.88%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'0.32%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'8.250"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-5.38%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1373"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.14%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.2412"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-0.03%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04410"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.05%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001257"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-0.33%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004308"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-7.32%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001201"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'6.25%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02789"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'11.53%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05398"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.23%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007787"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'4.31%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1415"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'0.76%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.008945"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'-7.22%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002172"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'2.62%"
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'11.05%"
$ws.Range("E45").Style = "Normal"
$ws.Range("E46").Value = "'-1.99%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Value = "'0.003203"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'7.80%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.002282"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.38%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002102"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
